$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Remove the two obsolete task rows ("個人社群網路建置" / "頭像編輯") that
#    used to sit at rows 21-22. Deleting them shifts every row below up by
#    two, which also removes the need to separately delete the two blank
#    trailing rows (60-61) - the sheet naturally ends up with 59 rows.
# ---------------------------------------------------------------------------
$ws.Rows("21:22").Delete()

# ---------------------------------------------------------------------------
# 2. Row 12 ("系統架構" task) - fill in the work-day count, the owner
#    ("全體") and the start/end dates (previously placeholder "?" cells).
# ---------------------------------------------------------------------------
$ws.Range("C12").Value = 1
$ws.Range("D19").Copy()
$ws.Range("D12").PasteSpecial(-4122)   # xlPasteFormats - copy D19's cell style (CJK font)
$ws.Range("D12").Value = "全體"
$ws.Range("E16").Copy()
$ws.Range("E12").PasteSpecial(-4122)   # copy E16's date number-format
$ws.Range("E12").Value = 40109
$ws.Range("F16").Copy()
$ws.Range("F12").PasteSpecial(-4122)
$ws.Range("F12").Value = 40109

# ---------------------------------------------------------------------------
# 3. Row 13 becomes the "編程標準" task, owned by "實作小組".
# ---------------------------------------------------------------------------
$ws.Range("B13").Value = "編程標準"
$ws.Range("C13").Value = 1
$ws.Range("D19").Copy()
$ws.Range("D13").PasteSpecial(-4122)
$ws.Range("D13").Value = "實作小組"
$ws.Range("E16").Copy()
$ws.Range("E13").PasteSpecial(-4122)
$ws.Range("E13").Value = 40110
$ws.Range("F16").Copy()
$ws.Range("F13").PasteSpecial(-4122)
$ws.Range("F13").Value = 40110

# ---------------------------------------------------------------------------
# 4. Row 14 becomes the "資料模型" task, also owned by "實作小組".
# ---------------------------------------------------------------------------
$ws.Range("B14").Value = "資料模型"
$ws.Range("C14").Value = 1
$ws.Range("D19").Copy()
$ws.Range("D14").PasteSpecial(-4122)
$ws.Range("D14").Value = "實作小組"
$ws.Range("E16").Copy()
$ws.Range("E14").PasteSpecial(-4122)
$ws.Range("E14").Value = 40110
$ws.Range("F16").Copy()
$ws.Range("F14").PasteSpecial(-4122)
$ws.Range("F14").Value = 40110

# ---------------------------------------------------------------------------
# 5. Row 22 (now "開發小組測試", shifted up from the old row 24) - finish
#    filling in its start date; the end date/owner were already correct.
# ---------------------------------------------------------------------------
$ws.Range("F22").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("E22").Value = 40113

# ---------------------------------------------------------------------------
# 6. Selection moves to E22, matching the saved view state.
# ---------------------------------------------------------------------------
$ws.Range("E22").Select()
